$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: collapse header row into a single title cell ---
$ws.Range("A1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# --- Rows 2-10: rearrange columns, per row mapping:
#   old A -> A (Day)
#   old B -> B (Time)
#   old E (Hours, numeric) -> C
#   old C (Module Code)    -> D
#   old D (Module Title)   -> E
#   old F -> F (Class Type)
#   old G -> G (Lecturer)
#   old J (Group)          -> H
#   old I -> I (Block / WLV)
#   old H (Room)            -> J
#   old K, old L            -> removed
for ($r = 2; $r -le 10; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2   # Day
    $b = $ws.Cells.Item($r, 2).Value2   # Time
    $c = $ws.Cells.Item($r, 3).Value2   # Module Code
    $d = $ws.Cells.Item($r, 4).Value2   # Module Title
    $e = $ws.Cells.Item($r, 5).Value2   # Hours
    $f = $ws.Cells.Item($r, 6).Value2   # Class Type
    $g = $ws.Cells.Item($r, 7).Value2   # Lecturer
    $h = $ws.Cells.Item($r, 8).Value2   # Room
    $i = $ws.Cells.Item($r, 9).Value2   # Block
    $j = $ws.Cells.Item($r, 10).Value2  # Group

    # Clear the whole row first so stale K/L (and anything else) is gone
    $rowRange = "A" + $r + ":L" + $r
    $ws.Range($rowRange).ClearContents()

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $e
    $ws.Cells.Item($r, 4).Value = $c
    $ws.Cells.Item($r, 5).Value = $d
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $j
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $h
}
